$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to carry a two-row header (units split across row 1 / row 2).
# Collapse that into a single header row by dropping the old row 2 - this
# shifts every data row up by one (old row 3..10 -> new row 2..9).
$ws.Rows.Item(2).Delete()

# Wipe whatever fragments of the old header remain on row 1 before writing
# the new, fuller header.
$ws.Range("A1:K1").ClearContents()

$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# Unit-row headers keep the small Arial font used elsewhere in the sheet.
$ws.Range("F1:K1").Font.Name = "Arial"
$ws.Range("F1:K1").Font.Size = 9

$ws.Range("A2:K2").Select()
